$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rows where the student marked Milestone "I" (roman numeral) in column E
# and "X" (completed) in column F.
$rows = @(4, 7, 21, 28, 37, 38)
foreach ($r in $rows) {
    $ws.Range("E$r").Value = "I"
    $ws.Range("F$r").Value = "X"
}

# Extra credit / bonus row: mark column C (Milestone I bonus) with "X"
$ws.Range("C91").Value = "X"

# Update the selected cell to reflect where the user left off editing
[void]$ws.Range("F14").Select()
